$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "gb_lab_construct" -> "lab_construct" (exclude modified sequences from alignment)
$ws.Range("E1").Value = "lab_construct"

# Update selection to E1 as recorded in the saved view state
$ws.Range("E1").Select()
